$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dialogue")

# Removed debug suffix from the (English) app name string
$ws.Range("B5").Value = "Shooting Stars"

# Leave the view scrolled/selected on the cell that was just edited
$ws.Activate()
$ws.Range("B5").Select()
